# Rename a set of PowerFactory short-circuit result column headers
# (shared-string table entries used as row-1 header labels on every
# worksheet) to match the corrected naming used by
# get_pf_shortcircuit_results.py / test_all_faults_and_cases.py:
#   *_deg            -> *_degree
#   *_bus_deg        -> *_degree   (redundant "bus" dropped too)
#   *_bus_pu         -> *_pu       (redundant "bus" dropped)
#   pf_q_*_from/to_mw -> pf_q_*_from/to_mvar (unit fix: MW -> Mvar)
#
# The header text only ever lives in row 1 of each worksheet, so a
# whole-cell, case-sensitive Find/Replace over every sheet's UsedRange
# is sufficient and avoids any accidental partial-text matches.

$wb = $excel.ActiveWorkbook

$renames = @(
    @{Old="pf_ikss_from_deg"; New="pf_ikss_from_degree"},
    @{Old="pf_ikss_to_deg"; New="pf_ikss_to_degree"},
    @{Old="pf_va_from_deg"; New="pf_va_from_degree"},
    @{Old="pf_va_to_deg"; New="pf_va_to_degree"},
    @{Old="pf_q_a_from_mw"; New="pf_q_a_from_mvar"},
    @{Old="pf_q_b_from_mw"; New="pf_q_b_from_mvar"},
    @{Old="pf_q_c_from_mw"; New="pf_q_c_from_mvar"},
    @{Old="pf_q_a_to_mw"; New="pf_q_a_to_mvar"},
    @{Old="pf_q_b_to_mw"; New="pf_q_b_to_mvar"},
    @{Old="pf_q_c_to_mw"; New="pf_q_c_to_mvar"},
    @{Old="pf_ikss_a_from_deg"; New="pf_ikss_a_from_degree"},
    @{Old="pf_ikss_b_from_deg"; New="pf_ikss_b_from_degree"},
    @{Old="pf_ikss_c_from_deg"; New="pf_ikss_c_from_degree"},
    @{Old="pf_ikss_a_to_deg"; New="pf_ikss_a_to_degree"},
    @{Old="pf_ikss_b_to_deg"; New="pf_ikss_b_to_degree"},
    @{Old="pf_ikss_c_to_deg"; New="pf_ikss_c_to_degree"},
    @{Old="pf_vm_b_from_bus_pu"; New="pf_vm_b_from_pu"},
    @{Old="pf_vm_c_from_bus_pu"; New="pf_vm_c_from_pu"},
    @{Old="pf_vm_a_to_bus_pu"; New="pf_vm_a_to_pu"},
    @{Old="pf_vm_b_to_bus_pu"; New="pf_vm_b_to_pu"},
    @{Old="pf_vm_c_to_bus_pu"; New="pf_vm_c_to_pu"},
    @{Old="pf_va_a_from_bus_deg"; New="pf_va_a_from_degree"},
    @{Old="pf_va_b_from_bus_deg"; New="pf_va_b_from_degree"},
    @{Old="pf_va_c_from_bus_deg"; New="pf_va_c_from_degree"},
    @{Old="pf_va_a_to_bus_deg"; New="pf_va_a_to_degree"},
    @{Old="pf_va_b_to_bus_deg"; New="pf_va_b_to_degree"},
    @{Old="pf_va_c_to_bus_deg"; New="pf_va_c_to_degree"}
)

foreach ($ws in $wb.Worksheets) {
    $headerRow = $ws.Range("A1:AQ1")
    foreach ($pair in $renames) {
        $headerRow.Replace($pair.Old, $pair.New, 1, 1, $false, $false, $false)
    }
}
